# Update the date strings in column B (rows 2-27) of the active sheet.
# These cells hold plain text (not real Excel dates), e.g. "04/10/2019".
# The commit rolls each date forward to a new value in July 2022.
#
# NOTE: Simply assigning a date-look-alike string via Range.Value / Value2 /
# Formula causes Excel's normal data-entry autodetection to convert the text
# into a real date serial number (and stamp a date NumberFormat on the
# cell) - exactly like typing such a string into a fresh "General" cell in
# real Excel. To keep these as plain text (matching the original file's
# shared-string cells with no explicit number format), we build each value
# as a text formula (forcing a string result), then copy/paste-special as
# values only. That collapses the formula down to a literal shared-string
# value without re-triggering date inference and without touching the
# cell's style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDates = @(
    "05/07/2022",
    "06/07/2022",
    "07/07/2022",
    "08/07/2022",
    "09/07/2022",
    "10/07/2022",
    "11/07/2022",
    "12/07/2022",
    "13/07/2022",
    "14/07/2022",
    "15/07/2022",
    "16/07/2022",
    "17/07/2022",
    "18/07/2022",
    "19/07/2022",
    "20/07/2022",
    "21/07/2022",
    "22/07/2022",
    "23/07/2022",
    "24/07/2022",
    "25/07/2022",
    "26/07/2022",
    "27/07/2022",
    "28/07/2022",
    "29/07/2022",
    "30/07/2022"
)

for ($i = 0; $i -lt $newDates.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 2)

    # Force the literal text into the cell via a string formula ...
    $cell.Formula = '="' + $newDates[$i] + '"'

    # ... then flatten the formula down to its plain text value so the
    # saved cell is a normal shared string again (no formula, no style
    # change, no date reinterpretation).
    $cell.Copy($cell)
    $cell.PasteSpecial(-4163)
}

$excel.CutCopyMode = $false
